$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Change 1: "Piyush Kumar Mallick 23079409" -> split into several runs with
# proofErr gramStart/gramEnd markers bracketing "Mallick  23079409"
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Piyush Kumar Mallick 23079409", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target1 = $d.Range($r1.Start, $r1.End)

$body1 = '<w:body><w:p>' +
  '<w:r w:rsidR="00852E4E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Piyush Kumar </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Mallick </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>23079409</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p></w:body>'

$target1.InsertXML($pkgOpen + $body1 + $pkgClose)

# ---------------------------------------------------------------------------
# Change 2: after "Ahmad Mujtaba Khan " insert an extra run with two spaces
# before the "23111308" run. (Insertion points degenerate to an empty
# range create a stray new paragraph, so instead we replace the
# non-degenerate "23111308" run with a two-space run followed by the same
# digits, split into two runs.)
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("23111308", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target2 = $d.Range($r2.Start, $r2.End)

$body2 = '<w:body><w:p>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r>' +
  '<w:r w:rsidR="00A96AFD" w:rsidRPr="00A96AFD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>23111308</w:t></w:r>' +
  '</w:p></w:body>'

$target2.InsertXML($pkgOpen + $body2 + $pkgClose)

# ---------------------------------------------------------------------------
# Change 3: insert an extra run with two spaces right before "22032491".
# Same technique as change 2.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("22032491", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target3 = $d.Range($r3.Start, $r3.End)

$body3 = '<w:body><w:p>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r>' +
  '<w:r w:rsidR="007D17EF" w:rsidRPr="007D17EF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>22032491</w:t></w:r>' +
  '</w:p></w:body>'

$target3.InsertXML($pkgOpen + $body3 + $pkgClose)

# ---------------------------------------------------------------------------
# Change 4: replace the two "[Name and ID ...]" placeholder paragraphs with a
# single paragraph containing "Haider Abid" / "23081929".
# ---------------------------------------------------------------------------
$startPara = $d.Content
$startPara.Find.Execute("[Name and ID of submitting student", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$endPara = $d.Content
$endPara.Find.Execute("[Name and ID of other group members]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target4 = $d.Range($startPara.Start, $endPara.End)

$body4 = '<w:body><w:p w14:paraId="61AE78EE" w14:textId="2546CBD5" w:rsidR="000F7A9B" w:rsidRPr="000F7A9B" w:rsidRDefault="000F7A9B">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:tab/><w:t xml:space="preserve">         Haider Abid                  </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>23081929</w:t></w:r>' +
  '</w:p></w:body>'

$target4.InsertXML($pkgOpen + $body4 + $pkgClose)

Write-Host "edits applied"
